# Refactor synthetic array /3: swap colored-square emoji for colored-book
# emoji in the "statut" column, and rename the "noir" (black) status label
# to "bleu" (blue) to match the new 📘 icon.
$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Cells.Replace("🟥", "📕")
$ws.Cells.Replace("⬛", "📘")
$ws.Cells.Replace("🟧", "📙")
$ws.Cells.Replace("🟩", "📗")
$ws.Cells.Replace("noir", "bleu")
